$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping CSV2openEHR")

$ws.Range("B7").Value = "bericht/altersangaben<<index>>/subject|id_namespace"
$ws.Range("C7").Value = 0
[void]$ws.Range("D11").Select()
